$wb = $excel.ActiveWorkbook

# Sheet "展览" (sheet1)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Cells.Item(2, 6).Value = 4739
$ws1.Cells.Item(3, 6).Value = 1894
$ws1.Cells.Item(6, 6).Value = 3204
$ws1.Cells.Item(10, 6).Value = 661
$ws1.Cells.Item(11, 6).Value = 557
$ws1.Cells.Item(12, 6).Value = 563
$ws1.Cells.Item(13, 6).Value = 416
$ws1.Cells.Item(14, 6).Value = 146
$ws1.Cells.Item(15, 6).Value = 1805
$ws1.Cells.Item(16, 6).Value = 1397
$ws1.Cells.Item(17, 6).Value = 129
$ws1.Cells.Item(18, 6).Value = 1660
$ws1.Cells.Item(20, 6).Value = 133
$ws1.Cells.Item(22, 6).Value = 21
$ws1.Cells.Item(26, 6).Value = 62
$ws1.Cells.Item(27, 6).Value = 120
$ws1.Cells.Item(30, 6).Value = 48
$ws1.Cells.Item(32, 6).Value = 4067
$ws1.Cells.Item(34, 6).Value = 786
$ws1.Cells.Item(36, 6).Value = 1799
$ws1.Cells.Item(38, 6).Value = 1927

# Sheet "演出" (sheet2)
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Cells.Item(2, 6).Value = 29
$ws2.Cells.Item(3, 6).Value = 59

# Sheet "全部类型" (sheet4)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Cells.Item(2, 6).Value = 4739
$ws4.Cells.Item(3, 6).Value = 1894
$ws4.Cells.Item(6, 6).Value = 3204
$ws4.Cells.Item(10, 6).Value = 661
$ws4.Cells.Item(11, 6).Value = 557
$ws4.Cells.Item(12, 6).Value = 563
$ws4.Cells.Item(13, 6).Value = 29
$ws4.Cells.Item(14, 6).Value = 416
$ws4.Cells.Item(15, 6).Value = 146
$ws4.Cells.Item(16, 6).Value = 1805
$ws4.Cells.Item(17, 6).Value = 1397
$ws4.Cells.Item(18, 6).Value = 129
$ws4.Cells.Item(19, 6).Value = 1660
$ws4.Cells.Item(21, 6).Value = 133
$ws4.Cells.Item(23, 6).Value = 21
$ws4.Cells.Item(27, 6).Value = 62
$ws4.Cells.Item(28, 6).Value = 120
$ws4.Cells.Item(31, 6).Value = 48
$ws4.Cells.Item(33, 6).Value = 4067
$ws4.Cells.Item(34, 6).Value = 59
$ws4.Cells.Item(37, 6).Value = 786
$ws4.Cells.Item(39, 6).Value = 1799
$ws4.Cells.Item(41, 6).Value = 1927
